$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Locate the "Micro results" row by checking the first cell's label text.
$targetRow = $null
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $label = $tbl.Rows.Item($i).Cells.Item(1).Range.Text
    if ($label -like "Micro results*") {
        $targetRow = $tbl.Rows.Item($i)
        break
    }
}

$cell = $targetRow.Cells.Item(2)

# Replace the single existing paragraph (with its paragraph mark) with a
# clean empty paragraph - this avoids leaving a stray <w:pPr>/<w:rPr> behind.
$firstParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/></w:p>
'@
$firstPara = $cell.Range.Paragraphs.Item(1)
$firstPara.Range.InsertXML($firstParaXml)

# Append the remaining new paragraphs (each its own run, blue 10pt Times New
# Roman) at the end of the cell.
$restParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – BLC – PERIPHERAL–RIGHT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – UC – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic recommendations without pathogen detection.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – CSF MICROSCOPY – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Gram stain shows no organisms.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – CSF CULT AND MICRO – CEREBROSPINAL FLUID;Brain NO GROWTH AFTER 2 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>13/06 – RESP. CULT AND MICRO – BRONCHO–ALVEOLAR LAVAGE NO GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>12/06 – CSF CULT AND MICRO – CEREBROSPINAL FLUIDHEAD NO GROWTH AFTER 10 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>12/06 – SARS–CoV–2 RNA – Negative</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>12/06 – SARS CORONAVIRUS–2 PCR – **No clear Result** +</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: MS2 control CT 23 detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>12/06 – CPE SCREEN – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: No carbapenem–resistant Enterobacteriaceae.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>11/06 – UC – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic recommendations; interpret with urine dipstick.  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Reasoning:** The report lacks explicit microbiological findings (e.g., organism growth, pathogen identification). It focuses on antibiotic guidance and clinical interpretation, not definitive infection status.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>11/06 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>11/06 – MRSA SCREEN – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: No MRSA isolated.</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($restParaXml)

Write-Output ("Paragraphs in cell: " + $cell.Range.Paragraphs.Count)
Write-Output $cell.Range.Text
